$p = $ppt.ActivePresentation

# In the target deck a new "Lessons Learned" slide is inserted right after
# the "HTML" slide (index 5) and right before the existing "Conclusion"
# slide. Its shapes/placeholders/formatting are identical to "Conclusion"
# (same shape ids, same creationId GUIDs, same xfrm) -- i.e. it was
# produced by duplicating the "Conclusion" slide and retitling the copy.
$conclusion = $p.Slides.Item(6)
$dup = $conclusion.Duplicate()
$newSlide = $dup.Item(1)

# Move the duplicate in front of the original "Conclusion" slide, pushing
# "Conclusion" down to index 7 and "Project Sources" down to index 8.
$newSlide.MoveTo(6)

# Retitle the duplicate's title placeholder (currently still reads
# "Conclusion"); leave the already-empty body placeholder untouched.
for ($i = 1; $i -le $newSlide.Shapes.Count; $i++) {
    $shp = $newSlide.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText -and $shp.TextFrame.TextRange.Text -eq "Conclusion") {
        $shp.TextFrame.TextRange.Text = "Lessons Learned"
        break
    }
}
